$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 4902
$ws.Range("I12").Value = 5360.1055
$ws.Range("J12").Value = 550
$ws.Range("K12").Value = 5360.1055
$ws.Range("L12").Value = 550
$ws.Range("M12").Value = -5190.1055
$ws.Range("N12").Value = -890

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 15901.9
$ws.Range("I21").Value = 16779.889
$ws.Range("J21").Value = 8000
$ws.Range("K21").Value = 16779.889
$ws.Range("L21").Value = 8000
$ws.Range("M21").Value = -16311.889
$ws.Range("N21").Value = -8936

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H23").Value = 15901.9
$ws.Range("I23").Value = 16779.889
$ws.Range("J23").Value = 8000
$ws.Range("K23").Value = 16779.889
$ws.Range("L23").Value = 8000
$ws.Range("M23").Value = -16545.889
$ws.Range("N23").Value = -8468

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 5099
$ws.Range("I31").Value = 200
$ws.Range("J31").Value = 9998
$ws.Range("K31").Value = 600
$ws.Range("L31").Value = 29994
$ws.Range("M31").Value = -370
$ws.Range("N31").Value = -30454

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 660
$ws.Range("I34").Value = 660
$ws.Range("K34").Value = 660
$ws.Range("M34").Value = -457

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H36").Value = 660
$ws.Range("I36").Value = 660
$ws.Range("K36").Value = 660
$ws.Range("M36").Value = 55

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 382.26086
$ws.Range("I80").Value = 323.42856
$ws.Range("K80").Value = 970.28568
$ws.Range("M80").Value = 27.71432000000004

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 382.26086
$ws.Range("I83").Value = 323.42856
$ws.Range("K83").Value = 2910.85704
$ws.Range("M83").Value = 2081.14296

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1729.7059
$ws.Range("I100").Value = 1700.4166
$ws.Range("J100").Value = 1800
$ws.Range("K100").Value = 1700.4166
$ws.Range("L100").Value = 1800
$ws.Range("M100").Value = -1159.4166
$ws.Range("N100").Value = -2882

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 241093.64
$ws.Range("I132").Value = 259516.23
$ws.Range("J132").Value = 1600
$ws.Range("K132").Value = 778548.6900000001
$ws.Range("L132").Value = 4800
$ws.Range("M132").Value = -776018.6900000001
$ws.Range("N132").Value = -9860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12363.655
$ws.Range("I32").Value = 5923.898
$ws.Range("J32").Value = 38659.332
$ws.Range("K32").Value = 5923.898
$ws.Range("L32").Value = 38659.332
$ws.Range("M32").Value = -5636.898
$ws.Range("N32").Value = -39233.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 871.8333
$ws.Range("I20").Value = 595.3333
$ws.Range("J20").Value = 1424.8334
$ws.Range("K20").Value = 595.3333
$ws.Range("L20").Value = 1424.8334
$ws.Range("M20").Value = -348.3333
$ws.Range("N20").Value = -1918.8334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1727.5
$ws.Range("I105").Value = 1817.1428
$ws.Range("J105").Value = 1518.3334
$ws.Range("K105").Value = 1817.1428
$ws.Range("L105").Value = 1518.3334
$ws.Range("M105").Value = -70.14280000000008
$ws.Range("N105").Value = -5012.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 265.8684
$ws.Range("I22").Value = 214.31429
$ws.Range("J22").Value = 867.3333
$ws.Range("K22").Value = 214.31429
$ws.Range("L22").Value = 867.3333
$ws.Range("M22").Value = 135.68571
$ws.Range("N22").Value = -1567.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1437
$ws.Range("I31").Value = 1134
$ws.Range("J31").Value = 1653.4286
$ws.Range("K31").Value = 1134
$ws.Range("L31").Value = 1653.4286
$ws.Range("M31").Value = -839
$ws.Range("N31").Value = -2243.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1437
$ws.Range("I34").Value = 1134
$ws.Range("J34").Value = 1653.4286
$ws.Range("K34").Value = 1134
$ws.Range("L34").Value = 1653.4286
$ws.Range("M34").Value = -932
$ws.Range("N34").Value = -2057.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H45").Value = 22500
$ws.Range("I45").Value = 10000
$ws.Range("K45").Value = 10000
$ws.Range("M45").Value = -9407

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 3813.762
$ws.Range("I107").Value = 4311.9165
$ws.Range("J107").Value = 3149.5557
$ws.Range("K107").Value = 4311.9165
$ws.Range("L107").Value = 3149.5557
$ws.Range("M107").Value = -2391.9165
$ws.Range("N107").Value = -6989.5557

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 9259729
$ws.Range("I5").Value = 485.54544
$ws.Range("J5").Value = 15625458
$ws.Range("K5").Value = 1456.63632
$ws.Range("L5").Value = 46876374
$ws.Range("M5").Value = -1344.63632
$ws.Range("N5").Value = -46876598

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1113.1265
$ws.Range("I68").Value = 788.9761999999999
$ws.Range("J68").Value = 1415.6666
$ws.Range("K68").Value = 2366.9286
$ws.Range("L68").Value = 4246.9998
$ws.Range("M68").Value = -1555.9286
$ws.Range("N68").Value = -5868.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1113.1265
$ws.Range("I71").Value = 788.9761999999999
$ws.Range("J71").Value = 1415.6666
$ws.Range("K71").Value = 7100.7858
$ws.Range("L71").Value = 12740.9994
$ws.Range("M71").Value = -3044.7858
$ws.Range("N71").Value = -20852.9994

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1565242.6
$ws.Range("J131").Value = 1820207.2
$ws.Range("L131").Value = 5460621.6
$ws.Range("N131").Value = -5470701.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 9259729
$ws.Range("I135").Value = 485.54544
$ws.Range("J135").Value = 15625458
$ws.Range("K135").Value = 4369.90896
$ws.Range("L135").Value = 140629122
$ws.Range("M135").Value = -1834.90896
$ws.Range("N135").Value = -140634192

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 44489.46
$ws.Range("I137").Value = 3436
$ws.Range("J137").Value = 70147.875
$ws.Range("K137").Value = 10308
$ws.Range("L137").Value = 210443.625
$ws.Range("M137").Value = -5208
$ws.Range("N137").Value = -220643.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 2788.875
$ws.Range("I30").Value = 2788.875
$ws.Range("K30").Value = 2788.875
$ws.Range("M30").Value = -2680.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 230.16667
$ws.Range("I55").Value = 282
$ws.Range("J55").Value = 193.14285
$ws.Range("K55").Value = 282
$ws.Range("L55").Value = 193.14285
$ws.Range("M55").Value = -109
$ws.Range("N55").Value = -539.14285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3925
$ws.Range("I61").Value = 1900
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 1900
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -1698
$ws.Range("N61").Value = -10404

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 3925
$ws.Range("I113").Value = 1900
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 1900
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = 270
$ws.Range("N113").Value = -14340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 14983.333
$ws.Range("J51").Value = 16980
$ws.Range("L51").Value = 16980
$ws.Range("N51").Value = -18000

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9762
$ws.Range("I62").Value = 11100
$ws.Range("J62").Value = 9167.333000000001
$ws.Range("K62").Value = 11100
$ws.Range("L62").Value = 9167.333000000001
$ws.Range("M62").Value = -10476
$ws.Range("N62").Value = -10415.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 9762
$ws.Range("I65").Value = 11100
$ws.Range("J65").Value = 9167.333000000001
$ws.Range("K65").Value = 55500
$ws.Range("L65").Value = 45836.665
$ws.Range("M65").Value = -52380
$ws.Range("N65").Value = -52076.665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5585.5
$ws.Range("I81").Value = 2100
$ws.Range("K81").Value = 4200
$ws.Range("M81").Value = -3139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 5585.5
$ws.Range("I84").Value = 2100
$ws.Range("K84").Value = 21000
$ws.Range("M84").Value = -15696

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 828.3333
$ws.Range("I107").Value = 325.5
$ws.Range("J107").Value = 1230.6
$ws.Range("K107").Value = 976.5
$ws.Range("L107").Value = 3691.8
$ws.Range("M107").Value = 943.5
$ws.Range("N107").Value = -7531.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 924.9286
$ws.Range("I136").Value = 907.72
$ws.Range("J136").Value = 1068.3334
$ws.Range("K136").Value = 2723.16
$ws.Range("L136").Value = 3205.0002
$ws.Range("M136").Value = -173.15999999999985
$ws.Range("N136").Value = -8305.0002
